$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 2010
$ws.Range("A2").Select() | Out-Null
